$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1898346.2
$ws.Range("I38").Value = 4032383.5
$ws.Range("J38").Value = 1424.1111
$ws.Range("K38").Value = 12097150.5
$ws.Range("L38").Value = 4272.3333
$ws.Range("M38").Value = -12096778.5
$ws.Range("N38").Value = -5016.3333
$ws.Range("H112").Value = 1087.25
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1087.25
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 3261.75
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -5477.75

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 12000
$ws.Range("J8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("N8").Value = -12288
$ws.Range("H32").Value = 38919.62
$ws.Range("I32").Value = 6523.237
$ws.Range("J32").Value = 214785.72
$ws.Range("K32").Value = 6523.237
$ws.Range("L32").Value = 214785.72
$ws.Range("M32").Value = -6236.237
$ws.Range("N32").Value = -215359.72
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H64").Value = 35811.668
$ws.Range("J64").Value = 35811.668
$ws.Range("L64").Value = 35811.668
$ws.Range("N64").Value = -36307.668
$ws.Range("H67").Value = 35811.668
$ws.Range("J67").Value = 35811.668
$ws.Range("L67").Value = 35811.668
$ws.Range("N67").Value = -37527.668
$ws.Range("H92").Value = 13137.5
$ws.Range("J92").Value = 13137.5
$ws.Range("L92").Value = 13137.5
$ws.Range("N92").Value = -18129.5
$ws.Range("H110").Value = 34521410
$ws.Range("I110").Value = 37078476
$ws.Range("J110").Value = 999.5
$ws.Range("K110").Value = 37078476
$ws.Range("L110").Value = 999.5
$ws.Range("M110").Value = -37076431
$ws.Range("N110").Value = -5089.5

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 79459
$ws.Range("J70").Value = 79459
$ws.Range("L70").Value = 79459
$ws.Range("N70").Value = -80045
$ws.Range("H73").Value = 79459
$ws.Range("J73").Value = 79459
$ws.Range("L73").Value = 79459
$ws.Range("N73").Value = -81487
$ws.Range("H99").Value = 2067.2856
$ws.Range("J99").Value = 2172.75
$ws.Range("L99").Value = 2172.75
$ws.Range("N99").Value = -5168.75

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 3004
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H97").Value = 33000
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H124").Value = 32857.2
$ws.Range("J124").Value = 32857.2
$ws.Range("L124").Value = 32857.2
$ws.Range("N124").Value = -37767.2

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2349.1333
$ws.Range("J34").Value = 2506.2144
$ws.Range("L34").Value = 7518.6432
$ws.Range("N34").Value = -7686.6432
$ws.Range("H58").Value = 1135.2
$ws.Range("I58").Value = 1003.3333
$ws.Range("J58").Value = 1191.7142
$ws.Range("K58").Value = 3009.9999
$ws.Range("L58").Value = 3575.1426
$ws.Range("M58").Value = -2881.9999
$ws.Range("N58").Value = -3831.1426
$ws.Range("H64").Value = 1950
$ws.Range("I64").Value = 800
$ws.Range("J64").Value = 2333.3333
$ws.Range("K64").Value = 2400
$ws.Range("L64").Value = 6999.999899999999
$ws.Range("M64").Value = -2130
$ws.Range("N64").Value = -7539.999899999999
$ws.Range("H67").Value = 1950
$ws.Range("I67").Value = 800
$ws.Range("J67").Value = 2333.3333
$ws.Range("K67").Value = 2400
$ws.Range("L67").Value = 6999.999899999999
$ws.Range("M67").Value = -1464
$ws.Range("N67").Value = -8871.999899999999
$ws.Range("H81").Value = 2651.8333
$ws.Range("I81").Value = 1042.4
$ws.Range("J81").Value = 3801.4285
$ws.Range("K81").Value = 3127.2
$ws.Range("L81").Value = 11404.2855
$ws.Range("M81").Value = -2004.2
$ws.Range("N81").Value = -13650.2855
$ws.Range("H84").Value = 2651.8333
$ws.Range("I84").Value = 1042.4
$ws.Range("J84").Value = 3801.4285
$ws.Range("K84").Value = 9381.6
$ws.Range("L84").Value = 34212.8565
$ws.Range("M84").Value = -3765.6
$ws.Range("N84").Value = -45444.8565
$ws.Range("H94").Value = 4325
$ws.Range("I94").Value = 3100
$ws.Range("J94").Value = 8000
$ws.Range("K94").Value = 9300
$ws.Range("L94").Value = 24000
$ws.Range("M94").Value = -8624
$ws.Range("N94").Value = -25352
$ws.Range("H107").Value = 1980.5454
$ws.Range("J107").Value = 1987.3334
$ws.Range("L107").Value = 5962.0002
$ws.Range("N107").Value = -9802.0002
$ws.Range("H131").Value = 810.6900000000001
$ws.Range("J131").Value = 861.5730600000001
$ws.Range("L131").Value = 2584.71918
$ws.Range("N131").Value = -12664.71918

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H15").Value = 12000
$ws.Range("J15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("N15").Value = -12576
$ws.Range("H48").Value = 13132.333
$ws.Range("J48").Value = 13132.333
$ws.Range("L48").Value = 13132.333
$ws.Range("N48").Value = -14102.333
$ws.Range("H81").Value = 12000
$ws.Range("J81").Value = 12000
$ws.Range("L81").Value = 12000
$ws.Range("N81").Value = -13996
$ws.Range("H84").Value = 12000
$ws.Range("J84").Value = 12000
$ws.Range("L84").Value = 36000
$ws.Range("N84").Value = -45984

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 318
$ws.Range("I17").Value = 318
$ws.Range("K17").Value = 318
$ws.Range("M17").Value = -148
$ws.Range("H22").Value = 1231.1666
$ws.Range("I22").Value = 846.75
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 846.75
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -551.75
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1231.1666
$ws.Range("I27").Value = 846.75
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 846.75
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -739.75
$ws.Range("N27").Value = -2214
$ws.Range("H46").Value = 1985.3529
$ws.Range("I46").Value = 300.2857
$ws.Range("K46").Value = 300.2857
$ws.Range("M46").Value = -112.2857
$ws.Range("H80").Value = 9673.529
$ws.Range("J80").Value = 9673.529
$ws.Range("L80").Value = 9673.529
$ws.Range("N80").Value = -11919.529
$ws.Range("H81").Value = 10122.667
$ws.Range("J81").Value = 10122.667
$ws.Range("L81").Value = 10122.667
$ws.Range("N81").Value = -12118.667
$ws.Range("H83").Value = 9673.529
$ws.Range("J83").Value = 9673.529
$ws.Range("L83").Value = 29020.587
$ws.Range("N83").Value = -40252.587
$ws.Range("H84").Value = 10122.667
$ws.Range("J84").Value = 10122.667
$ws.Range("L84").Value = 30368.001
$ws.Range("N84").Value = -40352.001
$ws.Range("H119").Value = 41985
$ws.Range("J119").Value = 41985
$ws.Range("L119").Value = 41985
$ws.Range("N119").Value = -51661

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1779.25
$ws.Range("I126").Value = 1814.5454
$ws.Range("J126").Value = 1649.8334
$ws.Range("K126").Value = 5443.6362
$ws.Range("L126").Value = 4949.5002
$ws.Range("M126").Value = -2973.6362
$ws.Range("N126").Value = -9889.5002
